$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell B11 ("R40" rule-name label) is changed to the text "1".
# A leading apostrophe forces Excel to store the numeric-looking
# entry as text (keeps it a shared string) instead of converting
# it to a number.
$ws.Range("B11").Value = "'1"
